$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 93
$ws.Range("H93").Value = 29690
$ws.Range("J93").Value = 29690
$ws.Range("L93").Value = 29690
$ws.Range("N93").Value = -34682
# Row 112
$ws.Range("H112").Value = 1293.7778
$ws.Range("J112").Value = 1333.0233
$ws.Range("L112").Value = 3999.0699
$ws.Range("N112").Value = -6215.0699
# Row 137
$ws.Range("H137").Value = 673596.9399999999
$ws.Range("I137").Value = 1645021.1
$ws.Range("J137").Value = 2851.6904
$ws.Range("K137").Value = 4935063.300000001
$ws.Range("L137").Value = 8555.0712
$ws.Range("M137").Value = -4932513.300000001
$ws.Range("N137").Value = -13655.0712
# Row 138
$ws.Range("H138").Value = 2420.4827
$ws.Range("I138").Value = 1387.8823
$ws.Range("J138").Value = 3883.3333
$ws.Range("K138").Value = 4163.6469
$ws.Range("L138").Value = 11649.9999
$ws.Range("M138").Value = 976.3531000000003
$ws.Range("N138").Value = -21929.9999
# Row 139
$ws.Range("H139").Value = 45325
$ws.Range("J139").Value = 45325
$ws.Range("L139").Value = 45325
$ws.Range("N139").Value = -55605

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7097.467
$ws.Range("I32").Value = 7844.5625
$ws.Range("J32").Value = 5258.4614
$ws.Range("K32").Value = 7844.5625
$ws.Range("L32").Value = 5258.4614
$ws.Range("M32").Value = -7557.5625
$ws.Range("N32").Value = -5832.4614
# Row 61
$ws.Range("H61").Value = 2037.5555
$ws.Range("I61").Value = 2037.5555
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2037.5555
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1825.5555
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 255036.14
$ws.Range("I74").Value = 355984
$ws.Range("J74").Value = 2666.5
$ws.Range("K74").Value = 355984
$ws.Range("L74").Value = 2666.5
$ws.Range("M74").Value = -355110
$ws.Range("N74").Value = -4414.5
# Row 77
$ws.Range("H77").Value = 255036.14
$ws.Range("I77").Value = 355984
$ws.Range("J77").Value = 2666.5
$ws.Range("K77").Value = 1779920
$ws.Range("L77").Value = 13332.5
$ws.Range("M77").Value = -1775552
$ws.Range("N77").Value = -22068.5
# Row 103
$ws.Range("H103").Value = 34666.668
$ws.Range("J103").Value = 34666.668
$ws.Range("L103").Value = 34666.668
$ws.Range("N103").Value = -37010.668
# Row 128
$ws.Range("H128").Value = 41980
$ws.Range("J128").Value = 41980
$ws.Range("L128").Value = 41980
$ws.Range("N128").Value = -51940
# Row 132
$ws.Range("H132").Value = 2555.9473
$ws.Range("I132").Value = 1916.5
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 5749.5
$ws.Range("L132").Value = 17899.0005
$ws.Range("M132").Value = -3219.5
$ws.Range("N132").Value = -22959.0005
# Row 136
$ws.Range("H136").Value = 2037.5555
$ws.Range("I136").Value = 2037.5555
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6112.666499999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3562.666499999999
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139").Value = 45250
$ws.Range("J139").Value = 45250
$ws.Range("L139").Value = 45250
$ws.Range("N139").Value = -55530

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2111.1
$ws.Range("I86").Value = 2080
$ws.Range("J86").Value = 2183.6667
$ws.Range("K86").Value = 2080
$ws.Range("L86").Value = 2183.6667
$ws.Range("M86").Value = -957
$ws.Range("N86").Value = -4429.6667
# Row 89
$ws.Range("H89").Value = 2111.1
$ws.Range("I89").Value = 2080
$ws.Range("J89").Value = 2183.6667
$ws.Range("K89").Value = 10400
$ws.Range("L89").Value = 10918.3335
$ws.Range("M89").Value = -4784
$ws.Range("N89").Value = -22150.3335
# Row 137
$ws.Range("H137").Value = 52190.625
$ws.Range("J137").Value = 52190.625
$ws.Range("L137").Value = 52190.625
$ws.Range("N137").Value = -62390.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 262413.56
$ws.Range("I31").Value = 483494.47
$ws.Range("J31").Value = 4485.8335
$ws.Range("K31").Value = 483494.47
$ws.Range("L31").Value = 4485.8335
$ws.Range("M31").Value = -483199.47
$ws.Range("N31").Value = -5075.8335
# Row 34
$ws.Range("H34").Value = 262413.56
$ws.Range("I34").Value = 483494.47
$ws.Range("J34").Value = 4485.8335
$ws.Range("K34").Value = 483494.47
$ws.Range("L34").Value = 4485.8335
$ws.Range("M34").Value = -483292.47
$ws.Range("N34").Value = -4889.8335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2578.5
$ws.Range("I3").Value = 1843
$ws.Range("J3").Value = 4785
$ws.Range("K3").Value = 5529
$ws.Range("L3").Value = 14355
$ws.Range("M3").Value = -5417
$ws.Range("N3").Value = -14579
# Row 68
$ws.Range("H68").Value = 1243.7344
$ws.Range("I68").Value = 1029
$ws.Range("J68").Value = 1445.4546
$ws.Range("K68").Value = 3087
$ws.Range("L68").Value = 4336.3638
$ws.Range("M68").Value = -2276
$ws.Range("N68").Value = -5958.3638
# Row 71
$ws.Range("H71").Value = 1243.7344
$ws.Range("I71").Value = 1029
$ws.Range("J71").Value = 1445.4546
$ws.Range("K71").Value = 9261
$ws.Range("L71").Value = 13009.0914
$ws.Range("M71").Value = -5205
$ws.Range("N71").Value = -21121.0914
# Row 81
$ws.Range("H81").Value = 2515
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2515
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 7545
$ws.Range("N81").Value = -9791
$ws.Range("M81").ClearContents()
# Row 84
$ws.Range("H84").Value = 2515
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2515
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 22635
$ws.Range("N84").Value = -33867
$ws.Range("M84").ClearContents()
# Row 93
$ws.Range("H93").Value = 10020
$ws.Range("J93").Value = 10020
$ws.Range("L93").Value = 30060
$ws.Range("N93").Value = -33804
# Row 94
$ws.Range("H94").Value = 3073.8333
$ws.Range("I94").Value = 999
$ws.Range("J94").Value = 3488.8
$ws.Range("K94").Value = 2997
$ws.Range("L94").Value = 10466.4
$ws.Range("M94").Value = -2321
$ws.Range("N94").Value = -11818.4
# Row 96
$ws.Range("H96").Value = 250003230
$ws.Range("J96").Value = 4309.3335
$ws.Range("L96").Value = 12928.0005
$ws.Range("N96").Value = -17046.0005
# Row 97
$ws.Range("H97").Value = 1100
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3300
$ws.Range("N97").Value = -4292
$ws.Range("M97").ClearContents()
# Row 105
$ws.Range("H105").Value = 7999
$ws.Range("J105").Value = 7999
$ws.Range("L105").Value = 23997
$ws.Range("N105").Value = -29239
# Row 106
$ws.Range("H106").Value = 3000
$ws.Range("J106").Value = 3000
$ws.Range("L106").Value = 9000
$ws.Range("N106").Value = -10892
# Row 107
$ws.Range("H107").Value = 28375.945
$ws.Range("I107").Value = 381.66666
$ws.Range("J107").Value = 103960.5
$ws.Range("K107").Value = 1144.99998
$ws.Range("L107").Value = 311881.5
$ws.Range("M107").Value = 775.0000199999999
$ws.Range("N107").Value = -315721.5
# Row 112
$ws.Range("H112").Value = 5713.636
$ws.Range("I112").Value = 462.5
$ws.Range("K112").Value = 1387.5
$ws.Range("M112").Value = -279.5
# Row 113
$ws.Range("H113").Value = 1712911.9
$ws.Range("I113").Value = 596.22644
$ws.Range("J113").Value = 6250548.5
$ws.Range("K113").Value = 1788.67932
$ws.Range("L113").Value = 18751645.5
$ws.Range("M113").Value = 381.3206799999998
$ws.Range("N113").Value = -18755985.5
# Row 129
$ws.Range("H129").Value = 1643.1
$ws.Range("I129").Value = 1037.8572
$ws.Range("J129").Value = 3055.3333
$ws.Range("K129").Value = 3113.5716
$ws.Range("L129").Value = 9165.999899999999
$ws.Range("M129").Value = 1886.4284
$ws.Range("N129").Value = -19165.9999
# Row 132
$ws.Range("H132").Value = 5069.9
$ws.Range("I132").Value = 999.5
$ws.Range("J132").Value = 6087.5
$ws.Range("K132").Value = 8995.5
$ws.Range("L132").Value = 54787.5
$ws.Range("M132").Value = -6465.5
$ws.Range("N132").Value = -59847.5
# Row 140
$ws.Range("H140").Value = 917
$ws.Range("I140").Value = 500.4
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 1501.2
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 3678.8
$ws.Range("N140").Value = -19360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1892.537
$ws.Range("I132").Value = 1661.4807
$ws.Range("J132").Value = 7900
$ws.Range("K132").Value = 4984.4421
$ws.Range("L132").Value = 23700
$ws.Range("M132").Value = -2454.4421
$ws.Range("N132").Value = -28760
# Row 137
$ws.Range("H137").Value = 39180
$ws.Range("J137").Value = 48770
$ws.Range("L137").Value = 48770
$ws.Range("N137").Value = -58970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4543.048
$ws.Range("I132").Value = 3615.9167
$ws.Range("J132").Value = 5779.222
$ws.Range("K132").Value = 10847.7501
$ws.Range("L132").Value = 17337.666
$ws.Range("M132").Value = -8317.750100000001
$ws.Range("N132").Value = -22397.666
# Row 136
$ws.Range("H136").Value = 2618.6086
$ws.Range("I136").Value = 1046.7
$ws.Range("J136").Value = 5565.9375
$ws.Range("K136").Value = 3140.1
$ws.Range("L136").Value = 16697.8125
$ws.Range("M136").Value = -590.1000000000004
$ws.Range("N136").Value = -21797.8125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2574.5667
$ws.Range("I132").Value = 1329.3125
$ws.Range("K132").Value = 3987.9375
$ws.Range("M132").Value = -1457.9375
# Row 138
$ws.Range("H138").Value = 44624.668
$ws.Range("J138").Value = 44624.668
$ws.Range("L138").Value = 44624.668
$ws.Range("N138").Value = -54904.668
